$wb = $excel.ActiveWorkbook

# --- Step 1: insert a new "2022-Q4" sheet before "2022-Q2" ---
# Copy the "2022-Q2" sheet (so it inherits identical sheet-level formatting:
# sheetPr/pageMargins/sheetFormatPr/header row style) and place the copy
# immediately before it, then rename the copy.
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($sheetQ2)
$q4 = $wb.Worksheets.Item("2022-Q2 (2)")
$q4.Name = "2022-Q4"

# --- Step 2: populate "2022-Q4" with the quarterly holdings data ---
# Extend column A's bordered/bold style (copied from A2) down through A10
# so every index cell in the new rows matches the existing formatting.
$q4.Range("A2").Copy()
$q4.Range("A3:A10").PasteSpecial(-4122)
$q4.Application.CutCopyMode = $false

# Force columns B:G to store as plain text (matches source data, which keeps
# numeric-looking fields such as fund codes/percentages as literal strings)
# by pre-formatting as Text before the values are entered.
$q4.Range("B2:G10").NumberFormat = "@"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "006049"
$q4.Cells.Item(2, 3).Value = "恒越研究精选混合A/B"
$q4.Cells.Item(2, 4).Value = "4.23"
$q4.Cells.Item(2, 5).Value = "88.62"
$q4.Cells.Item(2, 6).Value = "4.13"
$q4.Cells.Item(2, 7).Value = "0.1747"
$q4.Cells.Item(2, 8).Value = 5

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "007192"
$q4.Cells.Item(3, 3).Value = "恒越研究精选混合C"
$q4.Cells.Item(3, 4).Value = "3.19"
$q4.Cells.Item(3, 5).Value = "88.62"
$q4.Cells.Item(3, 6).Value = "4.13"
$q4.Cells.Item(3, 7).Value = "0.1317"
$q4.Cells.Item(3, 8).Value = 5

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "004784"
$q4.Cells.Item(4, 3).Value = "招商稳健优选股票"
$q4.Cells.Item(4, 4).Value = "3.60"
$q4.Cells.Item(4, 5).Value = "90.24"
$q4.Cells.Item(4, 6).Value = "1.91"
$q4.Cells.Item(4, 7).Value = "0.0688"
$q4.Cells.Item(4, 8).Value = 10

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "673141"
$q4.Cells.Item(5, 3).Value = "西部利得景程灵活配置混合A"
$q4.Cells.Item(5, 4).Value = "1.42"
$q4.Cells.Item(5, 5).Value = "82.89"
$q4.Cells.Item(5, 6).Value = "3.83"
$q4.Cells.Item(5, 7).Value = "0.0544"
$q4.Cells.Item(5, 8).Value = 6

$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 2).Value = "161224"
$q4.Cells.Item(6, 3).Value = "国投瑞银新丝路灵活配置混合（LOF）"
$q4.Cells.Item(6, 4).Value = "0.84"
$q4.Cells.Item(6, 5).Value = "93.28"
$q4.Cells.Item(6, 6).Value = "4.99"
$q4.Cells.Item(6, 7).Value = "0.0419"
$q4.Cells.Item(6, 8).Value = 6

$q4.Cells.Item(7, 1).Value = 5
$q4.Cells.Item(7, 2).Value = "002863"
$q4.Cells.Item(7, 3).Value = "金信深圳成长灵活配置混合"
$q4.Cells.Item(7, 4).Value = "0.73"
$q4.Cells.Item(7, 5).Value = "91.56"
$q4.Cells.Item(7, 6).Value = "5.09"
$q4.Cells.Item(7, 7).Value = "0.0372"
$q4.Cells.Item(7, 8).Value = 6

$q4.Cells.Item(8, 1).Value = 6
$q4.Cells.Item(8, 2).Value = "673143"
$q4.Cells.Item(8, 3).Value = "西部利得景程灵活配置混合C"
$q4.Cells.Item(8, 4).Value = "0.62"
$q4.Cells.Item(8, 5).Value = "82.89"
$q4.Cells.Item(8, 6).Value = "3.83"
$q4.Cells.Item(8, 7).Value = "0.0237"
$q4.Cells.Item(8, 8).Value = 6

$q4.Cells.Item(9, 1).Value = 7
$q4.Cells.Item(9, 2).Value = "012019"
$q4.Cells.Item(9, 3).Value = "国投瑞银安泽混合A"
$q4.Cells.Item(9, 4).Value = "0.62"
$q4.Cells.Item(9, 5).Value = "31.81"
$q4.Cells.Item(9, 6).Value = "0.98"
$q4.Cells.Item(9, 7).Value = "0.0061"
$q4.Cells.Item(9, 8).Value = 8

$q4.Cells.Item(10, 1).Value = 8
$q4.Cells.Item(10, 2).Value = "012020"
$q4.Cells.Item(10, 3).Value = "国投瑞银安泽混合C"
$q4.Cells.Item(10, 4).Value = "0.11"
$q4.Cells.Item(10, 5).Value = "31.81"
$q4.Cells.Item(10, 6).Value = "0.98"
$q4.Cells.Item(10, 7).Value = "0.0011"
$q4.Cells.Item(10, 8).Value = 8

# Drop back to the default (unstyled) cell format now that the text values are
# committed, so only the header row / column A keep the bold bordered style.
$q4.Range("B2:G10").Style = "Normal"

# --- Step 3: update "总计" (summary) sheet for the new quarter ---
# Row 2 becomes the new 2022-Q4 totals; the former row-2 (2022-Q2) and row-3
# (2022-Q1) entries shift down one row to make room.
$total = $wb.Worksheets.Item("总计")

# Extend column A's style from A3 down into the newly added A4.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 1
$total.Cells.Item(4, 4).Value = 0.03

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 1
$total.Cells.Item(3, 4).Value = 0.03

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 0.54

# Restore the originally-selected tab ("2022-Q1"), which simply shifted one
# position to the right when "2022-Q4" was inserted.
$wb.Worksheets.Item("2022-Q1").Activate()

